$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 203.33333
$ws.Range("I11").Value = 203.33333
$ws.Range("K11").Value = 203.33333
$ws.Range("M11").Value = -63.33332999999999
$ws.Range("H80").Value = 3174.4
$ws.Range("I80").Value = 1249.2
$ws.Range("J80").Value = 5099.6
$ws.Range("K80").Value = 3747.6
$ws.Range("L80").Value = 15298.8
$ws.Range("M80").Value = -2749.6
$ws.Range("N80").Value = -17294.8
$ws.Range("H83").Value = 3174.4
$ws.Range("I83").Value = 1249.2
$ws.Range("J83").Value = 5099.6
$ws.Range("K83").Value = 11242.8
$ws.Range("L83").Value = 45896.4
$ws.Range("M83").Value = -6250.800000000001
$ws.Range("N83").Value = -55880.4
$ws.Range("H99").Value = 187.2
$ws.Range("I99").Value = 186.75
$ws.Range("K99").Value = 560.25
$ws.Range("M99").Value = 937.75
$ws.Range("H101").Value = 233.85715
$ws.Range("J101").Value = 499
$ws.Range("L101").Value = 1497
$ws.Range("N101").Value = -4741
$ws.Range("H112").Value = 7300
$ws.Range("J112").Value = 7300
$ws.Range("L112").Value = 21900
$ws.Range("N112").Value = -24116
$ws.Range("H138").Value = 14832
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 14832
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 44496
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -54776

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 50000
$ws.Range("J18").Value = 50000
$ws.Range("L18").Value = 50000
$ws.Range("N18").Value = -50644
$ws.Range("H32").Value = 21531.527
$ws.Range("I32").Value = 20005.4
$ws.Range("K32").Value = 20005.4
$ws.Range("M32").Value = -19718.4
$ws.Range("H61").Value = 680.2857
$ws.Range("I61").Value = 729.5
$ws.Range("J61").Value = 385
$ws.Range("K61").Value = 729.5
$ws.Range("L61").Value = 385
$ws.Range("M61").Value = -517.5
$ws.Range("N61").Value = -809
$ws.Range("H74").Value = 2788.5
$ws.Range("I74").Value = 2788.5
$ws.Range("K74").Value = 2788.5
$ws.Range("M74").Value = -1914.5
$ws.Range("H77").Value = 2788.5
$ws.Range("I77").Value = 2788.5
$ws.Range("K77").Value = 13942.5
$ws.Range("M77").Value = -9574.5
$ws.Range("H88").Value = 2545.5
$ws.Range("I88").Value = 1956.7142
$ws.Range("K88").Value = 1956.7142
$ws.Range("M88").Value = -1550.7142
$ws.Range("H91").Value = 2545.5
$ws.Range("I91").Value = 1956.7142
$ws.Range("K91").Value = 1956.7142
$ws.Range("M91").Value = -552.7141999999999
$ws.Range("H97").Value = 3043.7
$ws.Range("I97").Value = 523.8570999999999
$ws.Range("K97").Value = 523.8570999999999
$ws.Range("M97").Value = -27.85709999999995
$ws.Range("H110").Value = 2278.28
$ws.Range("I110").Value = 1215.5652
$ws.Range("K110").Value = 1215.5652
$ws.Range("M110").Value = 829.4348
$ws.Range("H132").Value = 1266.25
$ws.Range("I132").Value = 949.0345
$ws.Range("K132").Value = 2847.1035
$ws.Range("M132").Value = -317.1035000000002
$ws.Range("H136").Value = 680.2857
$ws.Range("I136").Value = 729.5
$ws.Range("J136").Value = 385
$ws.Range("K136").Value = 2188.5
$ws.Range("L136").Value = 1155
$ws.Range("M136").Value = 361.5
$ws.Range("N136").Value = -6255

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 25000
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 25000
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2926.3572
$ws.Range("I31").Value = 2248
$ws.Range("J31").Value = 3830.8333
$ws.Range("K31").Value = 2248
$ws.Range("L31").Value = 3830.8333
$ws.Range("M31").Value = -1953
$ws.Range("N31").Value = -4420.8333
$ws.Range("H34").Value = 2926.3572
$ws.Range("I34").Value = 2248
$ws.Range("J34").Value = 3830.8333
$ws.Range("K34").Value = 2248
$ws.Range("L34").Value = 3830.8333
$ws.Range("M34").Value = -2046
$ws.Range("N34").Value = -4234.8333
$ws.Range("H68").Value = 42295
$ws.Range("J68").Value = 42295
$ws.Range("L68").Value = 42295
$ws.Range("N68").Value = -43793
$ws.Range("H71").Value = 42295
$ws.Range("J71").Value = 42295
$ws.Range("L71").Value = 126885
$ws.Range("N71").Value = -134373
$ws.Range("H74").Value = 51314
$ws.Range("J74").Value = 51314
$ws.Range("L74").Value = 51314
$ws.Range("N74").Value = -53062
$ws.Range("H77").Value = 51314
$ws.Range("J77").Value = 51314
$ws.Range("L77").Value = 153942
$ws.Range("N77").Value = -162678
$ws.Range("H86").Value = 9096.5
$ws.Range("I86").Value = 8995.5
$ws.Range("J86").Value = 9197.5
$ws.Range("K86").Value = 8995.5
$ws.Range("L86").Value = 9197.5
$ws.Range("M86").Value = -7872.5
$ws.Range("N86").Value = -11443.5
$ws.Range("H89").Value = 9096.5
$ws.Range("I89").Value = 8995.5
$ws.Range("J89").Value = 9197.5
$ws.Range("K89").Value = 44977.5
$ws.Range("L89").Value = 45987.5
$ws.Range("M89").Value = -39361.5
$ws.Range("N89").Value = -57219.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 6500
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H39").Value = 12624.5
$ws.Range("J39").Value = 12624.5
$ws.Range("L39").Value = 37873.5
$ws.Range("N39").Value = -38461.5
$ws.Range("H55").Value = 13866.444
$ws.Range("J55").Value = 13866.444
$ws.Range("L55").Value = 41599.33199999999
$ws.Range("N55").Value = -41953.33199999999
$ws.Range("H81").Value = 4732.1113
$ws.Range("I81").Value = 2599.8333
$ws.Range("K81").Value = 7799.499899999999
$ws.Range("M81").Value = -6676.499899999999
$ws.Range("H84").Value = 4732.1113
$ws.Range("I84").Value = 2599.8333
$ws.Range("K84").Value = 23398.4997
$ws.Range("M84").Value = -17782.4997
$ws.Range("H121").Value = 2333
$ws.Range("I121").Value = 2333
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 6999
$ws.Range("L121").Value = 0
$ws.Range("M121").Value = -5689
$ws.Range("N121").ClearContents()
$ws.Range("H122").Value = 2270.125
$ws.Range("I122").Value = 699.5
$ws.Range("J122").Value = 2793.6667
$ws.Range("K122").Value = 6295.5
$ws.Range("L122").Value = 25143.0003
$ws.Range("M122").Value = -3845.5
$ws.Range("N122").Value = -30043.0003
$ws.Range("H129").Value = 2440.0454
$ws.Range("J129").Value = 3003.923
$ws.Range("L129").Value = 9011.769
$ws.Range("N129").Value = -19011.769
$ws.Range("H131").Value = 2253.5
$ws.Range("I131").Value = 1730
$ws.Range("J131").Value = 2777
$ws.Range("K131").Value = 5190
$ws.Range("L131").Value = 8331
$ws.Range("M131").Value = -150
$ws.Range("N131").Value = -18411

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2947.25
$ws.Range("I126").Value = 3216.8
$ws.Range("K126").Value = 9650.400000000001
$ws.Range("M126").Value = -7180.400000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2226.6667
$ws.Range("I16").Value = 589.5
$ws.Range("K16").Value = 589.5
$ws.Range("M16").Value = -419.5
$ws.Range("H22").Value = 11314.714
$ws.Range("I22").Value = 12934.286
$ws.Range("J22").Value = 9695.143
$ws.Range("K22").Value = 12934.286
$ws.Range("L22").Value = 9695.143
$ws.Range("M22").Value = -12639.286
$ws.Range("N22").Value = -10285.143
$ws.Range("H27").Value = 11314.714
$ws.Range("I27").Value = 12934.286
$ws.Range("J27").Value = 9695.143
$ws.Range("K27").Value = 12934.286
$ws.Range("L27").Value = 9695.143
$ws.Range("M27").Value = -12827.286
$ws.Range("N27").Value = -9909.143
$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 5000
$ws.Range("K40").Value = 5000
$ws.Range("M40").Value = -4864
$ws.Range("H46").Value = 785
$ws.Range("I46").Value = 785
$ws.Range("K46").Value = 785
$ws.Range("M46").Value = -597
$ws.Range("H132").Value = 3465.9
$ws.Range("I132").Value = 3266
$ws.Range("K132").Value = 9798
$ws.Range("M132").Value = -7268
$ws.Range("H136").Value = 4392.2
$ws.Range("I136").Value = 2997
$ws.Range("K136").Value = 8991
$ws.Range("M136").Value = -6441

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1997.5
$ws.Range("I96").Value = 1997.5
$ws.Range("K96").Value = 1997.5
$ws.Range("M96").Value = -624.5
$ws.Range("H126").Value = 1000
$ws.Range("I126").Value = 1000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -530
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 1463.8572
$ws.Range("I136").Value = 1207.8334
$ws.Range("K136").Value = 3623.5002
$ws.Range("M136").Value = -1073.5002
